$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.578.30"
$ws.Range("E2").Value = "  +4.69%  "

# Row 3
$ws.Range("D3").Value = "3.421.80"
$ws.Range("E3").Value = "  +6.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("D5").Value = "'576.50"
$ws.Range("E5").Value = "  +6.37%  "

# Row 6
$ws.Range("D6").Value = "'156.72"
$ws.Range("E6").Value = "  +6.52%  "

# Row 7
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("D8").Value = "3.427.29"
$ws.Range("E8").Value = "  +5.89%  "

# Row 9
$ws.Range("E9").Value = "  +0.56%  "

# Row 10
$ws.Range("D10").Value = "'7.56"
$ws.Range("E10").Value = "  +2.70%  "

# Row 11
$ws.Range("D11").Value = "'0.122"
$ws.Range("E11").Value = "  +7.56%  "

# Row 12
$ws.Range("D12").Value = "'0.437"
$ws.Range("E12").Value = "  +0.21%  "

# Row 13
$ws.Range("D13").Value = "4.021.43"
$ws.Range("E13").Value = "  +6.78%  "

# Row 14
$ws.Range("E14").Value = "  -0.86%  "

# Row 15
$ws.Range("D15").Value = "'0.0000186"
$ws.Range("E15").Value = "  +7.24%  "

# Row 16
$ws.Range("D16").Value = "'27.30"
$ws.Range("E16").Value = "  +4.45%  "

# Row 17
$ws.Range("D17").Value = "63.667.08"
$ws.Range("E17").Value = "  +5.03%  "

# Row 18
$ws.Range("D18").Value = "3.427.48"
$ws.Range("E18").Value = "  +6.52%  "

# Row 19
$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  +1.77%  "

# Row 20
$ws.Range("D20").Value = "'14.26"
$ws.Range("E20").Value = "  +6.81%  "

# Row 21
$ws.Range("D21").Value = "'8.48"
$ws.Range("E21").Value = "  +1.91%  "

# Row 22
$ws.Range("D22").Value = "'390.59"
$ws.Range("E22").Value = "  +3.71%  "

# Row 23
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.25%  "

# Row 24
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.539"
$ws.Range("E24").Value = "  +2.13%  "

# Row 25
$ws.Range("E25").Value = "  +3.13%  "

# Row 26
$ws.Range("D26").Value = "'0.0000108"
$ws.Range("E26").Value = "  +20.36%  "

# Row 27
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  +9.49%  "

# Row 28
$ws.Range("D28").Value = "'0.180"
$ws.Range("E28").Value = "  +5.14%  "

# Row 29
$ws.Range("E29").Value = "  +0.09%  "

# Row 30
$ws.Range("D30").Value = "'6.68"
$ws.Range("E30").Value = "  +8.18%  "

# Row 31
$ws.Range("E31").Value = "  +13.60%  "

# Row 32
$ws.Range("D32").Value = "'2.04"
$ws.Range("E32").Value = "  +6.73%  "

# Row 33
$ws.Range("D33").Value = "'5.78"
$ws.Range("E33").Value = "  +8.55%  "

# Row 34
$ws.Range("D34").Value = "'23.53"
$ws.Range("E34").Value = "  +4.15%  "

# Row 35
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  -0.09%  "

# Row 36
$ws.Range("D36").Value = "'6.84"
$ws.Range("E36").Value = "  +3.09%  "

# Row 37
$ws.Range("D37").Value = "'1.49"
$ws.Range("E37").Value = "  +7.23%  "

# Row 38
$ws.Range("D38").Value = "'158.37"

# Row 39
$ws.Range("D39").Value = "'28.13"
$ws.Range("E39").Value = "  +6.18%  "

# Row 40
$ws.Range("D40").Value = "'0.0780"
$ws.Range("E40").Value = "  +9.13%  "

# Row 41
$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  +9.59%  "

# Row 42
$ws.Range("D42").Value = "2.869.16"
$ws.Range("E42").Value = "  +2.23%  "

# Row 43
$ws.Range("D43").Value = "'0.0320"
$ws.Range("E43").Value = "  +1.58%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'41.91"
$ws.Range("E44").Value = "  +4.81%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.770"
$ws.Range("E45").Value = "  +6.32%  "

# Row 46
$ws.Range("D46").Value = "'4.37"
$ws.Range("E46").Value = "  +2.38%  "

# Row 47
$ws.Range("E47").Value = "  +10.00%  "

# Row 48
$ws.Range("D48").Value = "3.475.16"
$ws.Range("E48").Value = "  +6.78%  "

# Row 49
$ws.Range("D49").Value = "'22.57"
$ws.Range("E49").Value = "  +8.08%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.38"
$ws.Range("E50").Value = "  +2.98%  "

# Row 51
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'295.57"
$ws.Range("E51").Value = "  +10.56%  "
